$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 with the new date value
$ws.Range("B2").Value = "Lunes 24/06/2024"

# Remove rows 3 through 11 (previously held additional date entries)
$ws.Range("A3:B11").EntireRow.Delete()
